# "Business Exceptions" sheet, rows 16-18 held three "not a valid path"
# exception rows that the UiPath bot had inserted in the wrong order
# (magic9843, testps0324, magic435). The database-populating run fixed the
# ordering so it reads testps0324, magic435, magic9843 (row 19 /
# test234234 was already correct and is left untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Business Exceptions")

$ws.Range("A16").Value = 'The Conclusion Evidence Location path: FOLDER\Mainfolder\Remediation_or_Justification Evidence\1-09-2020\RandomFolder1\ for the change: CHANGES - SOX Audit Report for testps0324.txt_07.01.73.eml made on 1/9/2020 is not a valid path.'
$ws.Range("A17").Value = 'The Conclusion Evidence Location path: FOLDER\Mainfolder\Remediation_or_Justification Evidence\1-22-2020\RandomFolder2\ for the change: CHANGES - SOX Audit Report for magic435.txt_07.01.73.eml made on 1/22/2020 is not a valid path.'
$ws.Range("A18").Value = 'The Conclusion Evidence Location path: FOLDER\Mainfolder\Remediation_or_Justification Evidence\Invalid Folder\1-14-2020\CHR0000291924.pdf for the change: CHANGES - SOX Audit Report for magic9843.txt_07.01.73.eml made on 1/14/2020 is not a valid path.'

